# fix bug and update screen
#
# Rooms.xlsx / Sheet1 updates:
#  - D2, D3, D4 ("Phong dac biet" for P201/P202/P203) shrink from the
#    stale "SSG101, SSC102, SSC101" list down to just "SSG101".
#  - D48 (HB/304R) and D50 (HB/309R) no longer carry the bogus long
#    "GDG401, WMT201, ..." special-room list (cleared).
#  - D51 (HB/202R) keeps its existing "LAB101, LAB211" value.
#  - Column D is narrowed back down to fit the now much shorter text.
#  - The active selection / view moves down near the bottom of the list
#    (previously parked past the used range at D53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data fixes -----------------------------------------------------
$ws.Range("D2").Value = "SSG101"
$ws.Range("D3").Value = "SSG101"
$ws.Range("D4").Value = "SSG101"

$ws.Range("D48").ClearContents()
$ws.Range("D50").ClearContents()

# --- column width now that the long text is gone ---------------------
$ws.Columns.Item(4).ColumnWidth = 13.9

# --- update the view/selection ---------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("D50").Select()
